# Fix StudentsImport template header: accept "nis" (in addition to "nisn").
# The student-number column header in the student import template is
# renamed from "nisn" to "nis".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "nis"
